$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up existing rows ---
# Row 2: civil status corrected from "Civil stat#B" to "Civil stat#L"
$ws.Range("F2").Value = "Civil stat#L"

# Row 4: house number "3c" -> "3c!"
$ws.Range("K4").Value = "3c!"

# Row 9 (Anton Ackermann): birthdate becomes an (invalid/pseudonymized) text
# date, and civil status corrected to "Civil stat#M"
$ws.Range("D9").Value = "31.02.1965"
$ws.Range("F9").Value = "Civil stat#M"

# Row 10 (Anita Ackermann-Abegger): birthdate changes to a new date value,
# civil status corrected to "Civil stat#M"
$ws.Range("D10").Value = 26320
$ws.Range("F10").Value = "Civil stat#M"

# --- Add new row 11: a new person ("Ali Apmann") sharing the same address
#     as the Ackermann family (same NNSS-style id/ family logic) ---
# Copy formatting down from row 10 first so number formats (NNSS, date) match.
$ws.Range("A10").Copy($ws.Range("A11"))
$ws.Range("D10").Copy($ws.Range("D11"))

$ws.Range("A11").Value = 7560000000009
$ws.Range("B11").Value = "Apmann"
$ws.Range("C11").Value = "Ali"
$ws.Range("D11").Value = 30246
$ws.Range("E11").Value = "Sex#W"
$ws.Range("F11").Value = "Civil stat#L"
$ws.Range("G11").Value = "Antragssteller"
$ws.Range("H11").Value = 3
$ws.Range("I11").Value = "CHResidenc#L"
$ws.Range("J11").Value = "Ackerstrasse"
$ws.Range("K11").Value = 11
$ws.Range("L11").Value = 804500
$ws.Range("M11").Value = "Zürich"
$ws.Range("N11").Value = 77777
$ws.Range("O11").Value = 1
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 0
$ws.Range("R11").Value = 0
$ws.Range("S11").Value = 0
$ws.Range("T11").Value = -30

# --- Misc view state matching the saved workbook ---
$ws.Range("C8").Select()
